# Append two more copies of the existing 5 data rows (rows 11-15) to the
# bottom of the sheet, as new rows 16-20 and 21-25 respectively. The
# dimension will grow from A1:P15 to A1:P25 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceStartRow = 11
$sourceEndRow = 15
$numCols = 16  # columns A..P

# Capture the source block of values (rows 11-15, columns A-P) using
# Value2, which (unlike Value) reliably returns usable scalars in this
# runtime.
$sourceValues = @()
for ($r = $sourceStartRow; $r -le $sourceEndRow; $r++) {
    $rowValues = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $rowValues += , $ws.Cells.Item($r, $c).Value2
    }
    $sourceValues += , $rowValues
}

# Append the captured block twice: rows 16-20 and rows 21-25.
for ($copy = 1; $copy -le 2; $copy++) {
    $destStartRow = $sourceEndRow + (($copy - 1) * ($sourceEndRow - $sourceStartRow + 1)) + 1
    for ($i = 0; $i -lt $sourceValues.Count; $i++) {
        $destRow = $destStartRow + $i
        $rowValues = $sourceValues[$i]
        for ($c = 1; $c -le $numCols; $c++) {
            $destCell = $ws.Cells.Item($destRow, $c)
            # Force text interpretation so date-like strings (e.g.
            # "1975-10-03") are not auto-converted into real dates, then
            # reset the style so no extra number-format style sticks to
            # the cell (keeps it identical to the default/unstyled cells
            # used throughout the rest of the sheet).
            $destCell.NumberFormat = "@"
            $destCell.Value2 = $rowValues[$c - 1]
            $destCell.Style = "Normal"
        }
    }
}
